$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before D (new most-recent reporting period), shifting D:K -> E:L.
$ws.Columns("D:D").Insert()

# Copy number formatting from the (now-shifted) old column D, now at E, onto new column D.
$ws.Columns("E:E").Copy()
$ws.Columns("D:D").PasteSpecial(-4122)

# Populate new column D with the latest period values.
$ws.Cells.Item(7, 4).Value = 43465
$ws.Cells.Item(8, 4).Value = 1064200
$ws.Cells.Item(9, 4).Value = 807800
$ws.Cells.Item(10, 4).Value = 256400
$ws.Cells.Item(12, 4).Value = "NA"
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(14, 4).Value = 333500
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(17, 4).Value = 1427700
$ws.Cells.Item(18, 4).Value = -363500
$ws.Cells.Item(20, 4).Value = 6300
$ws.Cells.Item(21, 4).Value = -282700
$ws.Cells.Item(22, 4).Value = 32500
$ws.Cells.Item(23, 4).Value = -389800
$ws.Cells.Item(24, 4).Value = -100
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(26, 4).Value = -389700
$ws.Cells.Item(27, 4).Value = -389700
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(29, 4).Value = 15600
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(32, 4).Value = -6300
$ws.Cells.Item(33, 4).Value = -374100
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(35, 4).Value = -374100
$ws.Cells.Item(38, 4).Value = 43465
$ws.Cells.Item(41, 4).Value = 47200
$ws.Cells.Item(42, 4).Value = 0
$ws.Cells.Item(43, 4).Value = 206900
$ws.Cells.Item(44, 4).Value = 479000
$ws.Cells.Item(45, 4).Value = 32800
$ws.Cells.Item(46, 4).Value = 766000
$ws.Cells.Item(47, 4).Value = 45000
$ws.Cells.Item(48, 4).Value = 177400
$ws.Cells.Item(49, 4).Value = 828700
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(52, 4).Value = 12600
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(54, 4).Value = 1829700
$ws.Cells.Item(57, 4).Value = 143200
$ws.Cells.Item(58, 4).Value = 1200
$ws.Cells.Item(59, 4).Value = 92600
$ws.Cells.Item(60, 4).Value = 236900
$ws.Cells.Item(61, 4).Value = 517500
$ws.Cells.Item(62, 4).Value = 45100
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(66, 4).Value = 799500
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(72, 4).Value = 63700
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(76, 4).Value = 1030100
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(80, 4).Value = 43465
$ws.Cells.Item(81, 4).Value = -374100
$ws.Cells.Item(83, 4).Value = 74500
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(89, 4).Value = 2400
$ws.Cells.Item(91, 4).Value = -24000
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(94, 4).Value = -75400
$ws.Cells.Item(96, 4).Value = 0
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(100, 4).Value = 6500
$ws.Cells.Item(101, 4).Value = -1500
$ws.Cells.Item(102, 4).Value = -68000